$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 425 in the "Betarraga" data
# table, pushing the existing rows 425..544 down to 426..545.
$ws.Rows.Item(425).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A425").Value = 10
$ws.Range("B425").Value = "Vega Modelo de Temuco"
$ws.Range("C425").Value = "La Araucanía"
$ws.Range("D425").Value = 44988
$ws.Range("E425").Value = 9
$ws.Range("F425").Value = 100114014
$ws.Range("G425").Value = "Betarraga"
$ws.Range("H425").Value = "Sin especificar"
$ws.Range("I425").Value = "Primera"
$ws.Range("J425").Value = 40
$ws.Range("K425").Value = 9000
$ws.Range("L425").Value = 9000
$ws.Range("M425").Value = 9000
$ws.Range("N425").Value = "$/docena de paquetes"
$ws.Range("O425").Value = "Provincia de Cautín"
$ws.Range("P425").Value = 750
$ws.Range("Q425").Value = 12
$ws.Range("R425").Value = "Hortaliza"
